# Applies the "m4-handguards" balance-pass edit described in the commit
# "i may have been too mean the first time around":
#  - tweak a bunch of horizontal/vertical recoil (E/F) and ergonomics (D)
#    values on existing handguards
#  - fix the 2A Armament Builder Series rows which had gotten shuffled out
#    of size order (7"/10"/12"/15")
#  - add new MK10RL bottom/side/top adapter rows and FAB Defense / NcStar
#    picatinny-rail-adapter rows at the bottom of the table
#  - widen column A a bit now that it holds more/longer internal names

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Simple recoil/ergonomics tweaks on existing rows (no row movement)
# ---------------------------------------------------------------------
$ws.Range("F4").Value = -16
$ws.Range("F5").Value = -15

$ws.Range("D6").Value = 0.32
$ws.Range("E6").Value = -10

$ws.Range("D7").Value = 0.36
$ws.Range("E7").Value = -12

$ws.Range("F20").Value = -8
$ws.Range("F21").Value = -10

$ws.Range("E28").Value = -5
$ws.Range("F28").Value = -9

$ws.Range("E30").Value = -8
$ws.Range("F30").Value = -10

$ws.Range("F36").Value = -8
$ws.Range("F37").Value = -9
$ws.Range("F38").Value = -12
$ws.Range("E39").Value = -15

$ws.Range("F45").Value = -10
$ws.Range("F46").Value = -13

$ws.Range("F50").Value = -20

# ---------------------------------------------------------------------
# 2. Re-sort the 2A Armament Builder Series Gen 1 rows back into size
#    order: 7" (31), 10" (32), 12" (33), 15" (34). Write final values
#    directly so the end state is correct regardless of starting order.
# ---------------------------------------------------------------------
$ws.Range("A31").Value = "2a_armament_builder_series_gen1_7inch_handguard"
$ws.Range("B31").Value = "2A Armament Builder Series Gen 1 7"""
$ws.Range("C31").Value = 15
$ws.Range("D31").Value = 0.2
$ws.Range("E31").Value = -6
$ws.Range("F31").Value = -5
$ws.Range("M31").Value = 700

$ws.Range("A32").Value = "2a_armament_builder_series_gen1_10inch_handguard"
$ws.Range("B32").Value = "2A Armament Builder Series Gen 1 10"""
$ws.Range("C32").Value = 13
$ws.Range("D32").Value = 0.25
$ws.Range("E32").Value = -8
$ws.Range("F32").Value = -7
$ws.Range("M32").Value = 750

$ws.Range("A33").Value = "2a_armament_builder_series_gen1_12inch_handguard"
$ws.Range("B33").Value = "2A Armament Builder Series Gen 1 12"""
$ws.Range("C33").Value = 11
$ws.Range("D33").Value = 0.28000000000000003
$ws.Range("E33").Value = -10
$ws.Range("F33").Value = -9
$ws.Range("M33").Value = 800

$ws.Range("A34").Value = "2a_armament_builder_series_gen1_15inch_handguard"
$ws.Range("B34").Value = "2A Armament Builder Series Gen 1 15"""
$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 0.33
$ws.Range("E34").Value = -12
$ws.Range("F34").Value = -11
$ws.Range("M34").Value = 900

# ---------------------------------------------------------------------
# 3. New rows 52-61: MK10RL adapters + FAB Defense / NcStar picatinny
#    rail adapters (row 51 and 57 are intentionally left blank spacers).
# ---------------------------------------------------------------------
$ws.Range("A52").Value = "mk10_bottom_adapter_large"
$ws.Range("B52").Value = "MK10RL Large"
$ws.Range("C52").Value = -1
$ws.Range("D52").Value = 0.1
$ws.Range("E52").Value = 2
$ws.Range("F52").Value = 2
$ws.Range("M52").Value = 0

$ws.Range("A53").Value = "mk10_bottom_adapter_small"
$ws.Range("B53").Value = "MK10RL Small"
$ws.Range("C53").Value = -1
$ws.Range("D53").Value = 0.04
$ws.Range("E53").Value = 2
$ws.Range("F53").Value = 2
$ws.Range("M53").Value = 0

$ws.Range("A54").Value = "mk10_side_adapter_large"
$ws.Range("B54").Value = "MK10RL Large"
$ws.Range("C54").Value = -1
$ws.Range("D54").Value = 0.1
$ws.Range("M54").Value = 0

$ws.Range("A55").Value = "mk10_side_adapter_small"
$ws.Range("B55").Value = "MK10RL Small"
$ws.Range("C55").Value = -1
$ws.Range("D55").Value = 0.04
$ws.Range("M55").Value = 0

$ws.Range("A56").Value = "mk10_top_adapter_large"
$ws.Range("B56").Value = "MK10RL Large"
$ws.Range("C56").Value = -1
$ws.Range("D56").Value = 0.1
$ws.Range("M56").Value = 0

$ws.Range("A58").Value = "fab_defense_upr_16_4_bottom"
$ws.Range("B58").Value = "FAB Defense UPR 16/4 Picatinny Rail"
$ws.Range("C58").Value = -3
$ws.Range("D58").Value = 0.06
$ws.Range("E58").Value = 2
$ws.Range("F58").Value = 2
$ws.Range("M58").Value = 400

$ws.Range("A59").Value = "ncstar_marsv2_bottom"
$ws.Range("B59").Value = "NcStar MARSV2 Picatinny Rail"
$ws.Range("C59").Value = -2
$ws.Range("D59").Value = 0.07
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 5
$ws.Range("M59").Value = 300

$ws.Range("A60").Value = "fab_defense_upr_16_4_top"
$ws.Range("B60").Value = "FAB Defense UPR 16/4 Picatinny Rail"
$ws.Range("C60").Value = -3
$ws.Range("D60").Value = 0.06
$ws.Range("M60").Value = 400

$ws.Range("A61").Value = "ncstar_marsv2_top"
$ws.Range("B61").Value = "NcStar MARSV2 Picatinny Rail"
$ws.Range("C61").Value = -2
$ws.Range("D61").Value = 0.07
$ws.Range("M61").Value = 300

# ---------------------------------------------------------------------
# 4. Formulas: extend the "N" (score) column formula down through the
#    newly added rows (51-61), and (re)materialise the formula for the
#    rows whose inputs changed above.
# ---------------------------------------------------------------------
$ws.Range("N3:N11").Formula = "=C3-D3*20-E3*0.8-F3*0.6-H3*5+I3*10+J3/300"
$ws.Range("N32:N39").Formula = "=C32-D32*20-E32*0.8-F32*0.6-H32*5+I32*10+J32/300"
$ws.Range("N51:N61").Formula = "=C51-D51*20-E51*0.8-F51*0.6-H51*5+I51*10+J51/300"

# N31 keeps its own (non-shared) formula, same shape as before.
$ws.Range("N31").Formula = "=C31-D31*20-E31*0.8-F31*0.6-H31*5+I31*10+J31/300"

# ---------------------------------------------------------------------
# 5. Column A is now wide enough to show the longer internal names.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.8

# ---------------------------------------------------------------------
# 6. Tidy up the view: drop the old scrolled-down position and leave
#    the selection on Q17 (matches the saved view state).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q17").Select()

$wb.Application.Calculate()
